$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "ContosoLearn Competitor SWOT" "ContosoLearn Mitbewerber SWOT"

# Section headers (bold "Strengths:", "Weaknesses:", "Opportunities:", "Threats:") used twice each
Replace-Text "Strengths:" "Stärken:"
Replace-Text "Weaknesses:" "Schwachstellen:"
Replace-Text "Opportunities:" "Verkaufschancen:"
Replace-Text "Threats:" "Bedrohungen:"

# Fabrikam Learning - Strengths body
Replace-Text " Fabrikam Learning provides a comprehensive set of analytics and reporting tools. It ensures the continuous monitoring of teaching and learning activities, as well as pinpointing problematic areas that need to be addressed." " Fabrikam Learning bietet ein umfassendes Set an Analyse- und Berichterstattungs-Tools. Es gewährleistet die kontinuierliche Überwachung von Lehr- und Lernaktivitäten sowie die Anheftung problematischer Bereiche, die angegangen werden müssen."

# Fabrikam Learning - Weaknesses body
Replace-Text " While Fabrikam Learning has robust reporting capabilities, it might be overwhelming for some users due to its comprehensive nature." " Während Fabrikam Learning robuste Berichterstellungsfunktionen bietet, kann es aufgrund seiner umfassenden Natur für einige Benutzende überwältigend sein."

# Fabrikam Learning - Opportunities body
Replace-Text " There is a growing demand for personalized learning experiences and data-driven recommendations. Fabrikam Learning can leverage its robust analytics and reporting tools to meet this demand." " Es besteht eine wachsende Nachfrage nach personalisierten Lernerfahrungen und datengesteuerten Empfehlungen. Fabrikam Learning kann seine robusten Analyse- und Berichterstellungstools nutzen, um diese Nachfrage zu erfüllen."

# Fabrikam Learning - Threats body
Replace-Text " The eLearning market is highly competitive with many players offering similar features. Fabrikam Learning needs to continuously innovate to stay ahead." " Der eLearning-Markt ist mit vielen Anbietern, die ähnliche Funktionen anbieten, sehr wettbewerbsfähig. Fabrikam Learning muss kontinuierlich innovativ sein, um an der Spitze zu bleiben."

# AdatumLearn - Strengths body
Replace-Text " AdatumLearn offers courses on business analysis techniques such as MOST and SWOT. This shows their commitment to providing valuable content to their users." " AdatumLearn bietet Kurse zu Geschäftsanalysetechniken wie MOST und SWOT an. Dies zeigt ihr Engagement für die Bereitstellung wertvoller Inhalte für ihre Benutzenden."

# AdatumLearn - Weaknesses body
Replace-Text " The information provided in their courses is a compilation of third-party generated information. This might not be as valuable as original content." " Die in ihren Kursen bereitgestellten Informationen sind eine Zusammenstellung von generierten Informationen von Drittanbietern. Dies ist möglicherweise nicht so wertvoll wie der ursprüngliche Inhalt."

# AdatumLearn - Opportunities body
Replace-Text " AdatumLearn can create more original content to provide unique value to their users. They can also expand their course offerings to cover more topics." " AdatumLearn kann originellere Inhalte erstellen, um ihren Benutzenden einen einzigartigen Wert zu bieten. Sie können auch ihre Kursangebote erweitern, um weitere Themen zu behandeln."

# AdatumLearn - Threats body (original had a stray trailing quote character to drop)
Replace-Text " Like Fabrikam Learning, AdatumLearn also faces stiff competition in the eLearning market. They need to continuously improve their offerings to stay competitive.`"" " Wie Fabrikam Learning steht AdatumLearn auch im eLearning-Markt vor einem harten Wettbewerb. Sie müssen ihr Angebot kontinuierlich verbessern, um wettbewerbsfähig zu bleiben."
